# Iraq League workbook update (03-06-2024 23:01)
# - Swap several pairs of rows (columns B:AD) because the source data was
#   re-sorted by match id within identical kickoff dates.
# - Append 6 new match rows (228-233) at the end of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-RowValues($sheet, $row, $colStart, $colEnd) {
    $vals = @{}
    for ($c = $colStart; $c -le $colEnd; $c++) {
        $vals[$c] = $sheet.Cells.Item($row, $c).Value2
    }
    return $vals
}

function Clear-RowRange($sheet, $row, $colStart, $colEnd) {
    for ($c = $colStart; $c -le $colEnd; $c++) {
        $sheet.Cells.Item($row, $c).Value2 = $null
    }
}

function Set-RowValues($sheet, $row, $colStart, $colEnd, $vals) {
    for ($c = $colStart; $c -le $colEnd; $c++) {
        $v = $vals[$c]
        if ($v -ne $null) {
            $sheet.Cells.Item($row, $c).Value2 = $v
        }
    }
}

function Swap-Rows($sheet, $rowA, $rowB, $colStart, $colEnd) {
    $valsA = Get-RowValues $sheet $rowA $colStart $colEnd
    $valsB = Get-RowValues $sheet $rowB $colStart $colEnd

    Clear-RowRange $sheet $rowA $colStart $colEnd
    Clear-RowRange $sheet $rowB $colStart $colEnd

    Set-RowValues $sheet $rowA $colStart $colEnd $valsB
    Set-RowValues $sheet $rowB $colStart $colEnd $valsA
}

# Columns B..AD => 2..30
$colStart = 2
$colEnd = 30

$pairs = @(
    @(22,23),
    @(41,42),
    @(54,55),
    @(91,92),
    @(103,104),
    @(135,136),
    @(190,191),
    @(219,220)
)

foreach ($p in $pairs) {
    Swap-Rows $ws $p[0] $p[1] $colStart $colEnd
}

# New rows appended after row 227
function Set-MatchRow($sheet, $row, $data) {
    $sheet.Cells.Item($row, 1).Value2 = $data.A
    $sheet.Cells.Item($row, 2).Value2 = $data.B
    $sheet.Cells.Item($row, 3).Value2 = $data.C
    $sheet.Cells.Item($row, 4).Value2 = $data.D
    $sheet.Cells.Item($row, 5).Value2 = $data.E
    $sheet.Cells.Item($row, 6).Value2 = $data.F
    $sheet.Cells.Item($row, 7).Value2 = $data.G
    $sheet.Cells.Item($row, 8).Value2 = $data.H
    $sheet.Cells.Item($row, 11).Value2 = $data.K
    $sheet.Cells.Item($row, 12).Value2 = $data.L
    $sheet.Cells.Item($row, 13).Value2 = $data.M
    $sheet.Cells.Item($row, 14).Value2 = $data.N
    $sheet.Cells.Item($row, 15).Value2 = $data.O
    $sheet.Cells.Item($row, 16).Value2 = $data.P
    $sheet.Cells.Item($row, 17).Value2 = $data.Q
    $sheet.Cells.Item($row, 18).Value2 = $data.R
    $sheet.Cells.Item($row, 19).Value2 = $data.S
    $sheet.Cells.Item($row, 20).Value2 = $data.T
    $sheet.Cells.Item($row, 21).Value2 = $data.U
    $sheet.Cells.Item($row, 22).Value2 = $data.V
    $sheet.Cells.Item($row, 23).Value2 = $data.W
    $sheet.Cells.Item($row, 24).Value2 = $data.X
    $sheet.Cells.Item($row, 25).Value2 = $data.Y
    $sheet.Cells.Item($row, 26).Value2 = $data.Z
    $sheet.Cells.Item($row, 27).Value2 = $data.AA
    $sheet.Cells.Item($row, 28).Value2 = $data.AB
    $sheet.Cells.Item($row, 29).Value2 = $data.AC
    $sheet.Cells.Item($row, 30).Value2 = $data.AD
}

$row228 = @{
    A = 226; B = 8277157; C = "Iraq League"; D = 45442.5;
    E = "Naft Maysan"; F = "Al Shorta SC"; G = 2; H = 3;
    K = "A";
    L = 5.5; M = 3.4; N = 1.571; O = 5.75; P = 3.25; Q = 1.6; R = 1;
    S = 1.85; T = 1.95; U = 2.25; V = 1.8; W = 2;
    X = -1; Y = -1; Z = 0.6000000000000001; AA = 0; AB = 0; AC = 0.8; AD = -1
}

$row229 = @{
    A = 227; B = 8277170; C = "Iraq League"; D = 45442.60416666666;
    E = "Al Talaba"; F = "Naft AlWasat"; G = 1; H = 2;
    K = "A";
    L = 1.5; M = 4; N = 5; O = 1.4; P = 4.5; Q = 6; R = -1.25;
    S = 1.925; T = 1.875; U = 2.5; V = 1.9; W = 1.9;
    X = -1; Y = -1; Z = 5; AA = -1; AB = 0.875; AC = 0.8999999999999999; AD = -1
}

$row230 = @{
    A = 228; B = 8279463; C = "Iraq League"; D = 45443.42708333334;
    E = "Newroz SC"; F = "Al Quwa Al Jawiya"; G = 4; H = 2;
    K = "H";
    L = 4; M = 3.1; N = 1.85; O = 5.25; P = 3.4; Q = 1.6; R = 0.75;
    S = 2; T = 1.8; U = 2.5; V = 1.95; W = 1.85;
    X = 4.25; Y = -1; Z = -1; AA = 1; AB = -1; AC = 0.95; AD = -1
}

$row231 = @{
    A = 229; B = 8279461; C = "Iraq League"; D = 45443.52083333334;
    E = "Al Naft SC"; F = "Karbalaa FC"; G = 1; H = 1;
    K = "D";
    L = 2.25; M = 3; N = 3; O = 1.95; P = 3.1; Q = 3.75; R = -0.5;
    S = 2; T = 1.8; U = 2; V = 1.925; W = 1.875;
    X = -1; Y = 2.1; Z = -1; AA = -1; AB = 0.8; AC = 0; AD = 0
}

$row232 = @{
    A = 230; B = 8285514; C = "Iraq League"; D = 45443.52083333334;
    E = "Al Najaf"; F = "Zakho"; G = 2; H = 1;
    K = "H";
    L = 2.1; M = 3; N = 3.25; O = 2.2; P = 2.8; Q = 3.3; R = -0.25;
    S = 1.95; T = 1.85; U = 1.75; V = 1.85; W = 1.95;
    X = 1.2; Y = -1; Z = -1; AA = 0.95; AB = -1; AC = 0.8500000000000001; AD = -1
}

$row233 = @{
    A = 231; B = 8279462; C = "Iraq League"; D = 45443.61458333334;
    E = "Al Zawraa"; F = "Al Minaa"; G = 1; H = 0;
    K = "H";
    L = 1.25; M = 5; N = 9; O = 1.3; P = 4.75; Q = 7.5; R = -1.5;
    S = 1.925; T = 1.875; U = 2.5; V = 1.95; W = 1.85;
    X = 0.3; Y = -1; Z = -1; AA = -1; AB = 0.875; AC = -1; AD = 0.8500000000000001
}

Set-MatchRow $ws 228 $row228
Set-MatchRow $ws 229 $row229
Set-MatchRow $ws 230 $row230
Set-MatchRow $ws 231 $row231
Set-MatchRow $ws 232 $row232
Set-MatchRow $ws 233 $row233

# Replicate the formatting (bold/bordered id column, date number format
# column, etc.) from the last pre-existing data row onto the new rows,
# without creating any additional cell-style definitions.
$ws.Range("A227:AD227").Copy()
$ws.Range("A228:AD233").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# The source row (227) has no HTHG/HTAG (I/J) values, so the format-only
# paste above leaves behind empty I/J placeholder cells on the new rows.
# None of the new matches have half-time scores either, so drop them.
$ws.Range("I228:J233").ClearContents()

"Done"
